$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Multi-threaded (E/F/G) table: fix currency-converted score for
#     "GLOBAL PebbleHost Extreme 6GB" row (row 3) ---
$ws.Range("F3").Value2 = 105
$ws.Range("G3").Value2 = 38.43

# --- Single-threaded (A/B/C) table: add the missing
#     "EU VolcanoHosting Premium 4GB ($9.51)" entry, shifting the two
#     rows below it down by one ---
# Before:
#   A11 NA PebbleHost Premium 4GB ($9.00)     B11 49
#   A12 EU Birdflop Premium 4GB ($8.00)       B12 176
#   A13 EU VolcanoHosting Premium 4GB ($8.00) B13 68
# After:
#   A11 EU VolcanoHosting Premium 4GB ($9.51) B11 68
#   A12 NA PebbleHost Premium 4GB ($9.00)     B12 49
#   A13 EU Birdflop Premium 4GB ($8.00)       B13 176
$ws.Range("A13").Value2 = $ws.Range("A12").Value2
$ws.Range("B13").Value2 = $ws.Range("B12").Value2
$ws.Range("A12").Value2 = $ws.Range("A11").Value2
$ws.Range("B12").Value2 = $ws.Range("B11").Value2
$ws.Range("A11").Value2 = "EU VolcanoHosting Premium 4GB (`$9.51)"
$ws.Range("B11").Value2 = 68

# --- Cosmetic: reflect the author's last active-cell selection ---
$ws.Range("V20").Select()

$wb.Save()
